# Auto-update predictions and index for 2025-10-22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Arsenal FC vs Atlético de Madrid ---
$ws.Range("A2").Value = "Arsenal FC ✓ - Atlético de Madrid: 4:0"
$ws.Range("E2").Value = 82
$ws.Range("G2").Value = "✓"

# --- Row 3: Union Saint-Gilloise vs Inter Milan ---
$ws.Range("A3").Value = "Union Saint-Gilloise - Inter Milan ✓: 0:4"
$ws.Range("C3").Value = 70
$ws.Range("D3").Value = 92
$ws.Range("G3").Value = "✓"

# --- Row 4: new fixture, FC Barcelona vs Olympiacos Piraeus ---
$ws.Range("A4").Value = "FC Barcelona ✓ - Olympiacos Piraeus: 6:1"
$ws.Range("B4").Value = "FC Barcelona"
$ws.Range("C4").Value = 67
$ws.Range("D4").Value = 92
$ws.Range("E4").Value = 93
$ws.Range("F4").Value = 1.28
$ws.Range("G4").Value = "✓"
